# Lecture 06 - Image Descriptors
# Remove the "Classification" bullet item from the "Topics" slide (slide 2),
# leaving "Discussion of Lecture #05", "Feature Vector", "Horizontal and
# Vertical Projections", "Image Descriptors" and "Practice" in place with
# their original formatting untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(4)
$textRange = $shape.TextFrame.TextRange

# Find the "Classification" paragraph (5th paragraph in the bullet list) and
# delete it entirely, including its trailing paragraph mark, so the
# following "Practice" paragraph keeps its own original formatting.
for ($i = $textRange.Paragraphs().Count; $i -ge 1; $i--) {
    $para = $textRange.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd("`r")
    if ($paraText -eq "Classification") {
        $para.Delete()
    }
}
